$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9736156463623047
$ws.Range("B1").Value = 1.440665364265442
$ws.Range("C1").Value = 5.798338890075684
$ws.Range("D1").Value = 1.697846651077271
$ws.Range("E1").Value = 1.045493960380554
